# Auto-generated script to update cryptos.xlsx price/volume columns
# per commit 'Updated symbol list on Tue Feb  7 15:36:28 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the affected cells so values like '328.98' or
# '0.61%' are stored as literal text (matching the existing inline-string cells)
# instead of being auto-converted to numbers/percentages by Excel's input parser.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "328.98"
$ws.Range("E2").Value = "0.61%"
$ws.Range("D3").Value = "44.17"
$ws.Range("E3").Value = "0.84%"
$ws.Range("D4").Value = "5.485"
$ws.Range("E4").Value = "-0.93%"
$ws.Range("D5").Value = "0.08047"
$ws.Range("E5").Value = "0.43%"
$ws.Range("D6").Value = "2.028"
$ws.Range("E6").Value = "6.92%"
$ws.Range("D7").Value = "0.9532"
$ws.Range("E7").Value = "0.91%"
$ws.Range("D8").Value = "0.1112"
$ws.Range("E8").Value = "-5.45%"
$ws.Range("D9").Value = "0.1867"
$ws.Range("E9").Value = "1.69%"
$ws.Range("D10").Value = "10.22"
$ws.Range("E10").Value = "1.44%"
$ws.Range("D11").Value = "0.09830"
$ws.Range("E11").Value = "1.81%"
$ws.Range("D12").Value = "0.04749"
$ws.Range("E12").Value = "6.95%"
$ws.Range("D13").Value = "0.1064"
$ws.Range("E13").Value = "-0.21%"
$ws.Range("D14").Value = "0.001275"
$ws.Range("E14").Value = "-0.77%"
$ws.Range("E15").Value = "-2.55%"
$ws.Range("D16").Value = "0.005897"
$ws.Range("E16").Value = "-0.81%"
$ws.Range("D17").Value = "3.377"
$ws.Range("E17").Value = "-0.86%"
$ws.Range("D18").Value = "4.414"
$ws.Range("E18").Value = "3.51%"
$ws.Range("E19").Value = "3.27%"
$ws.Range("D20").Value = "0.3408"
$ws.Range("E20").Value = "-1.03%"
$ws.Range("D21").Value = "0.1402"
$ws.Range("E21").Value = "-0.30%"
$ws.Range("E22").Value = "2.84%"
$ws.Range("D23").Value = "0.001310"
$ws.Range("E23").Value = "4.93%"
$ws.Range("D24").Value = "0.004354"
$ws.Range("E24").Value = "1.60%"
$ws.Range("D25").Value = "0.0001252"
$ws.Range("E25").Value = "-0.94%"
$ws.Range("D26").Value = "0.0003750"
$ws.Range("E26").Value = "-6.17%"
$ws.Range("D38").Value = "0.02574"
$ws.Range("E38").Value = "-1.98%"
$ws.Range("D39").Value = "0.05648"
$ws.Range("E39").Value = "3.15%"
$ws.Range("D40").Value = "0.007739"
$ws.Range("E40").Value = "2.08%"
$ws.Range("D41").Value = "0.1398"
$ws.Range("E41").Value = "0.69%"
$ws.Range("D42").Value = "0.007367"
$ws.Range("E42").Value = "-9.97%"
$ws.Range("D43").Value = "0.002011"
$ws.Range("E43").Value = "0.28%"
$ws.Range("D44").Value = "0.008533"
$ws.Range("E44").Value = "-3.06%"
$ws.Range("D45").Value = "0.00007111"
$ws.Range("E45").Value = "-0.08%"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").Value = "0.09%"
$ws.Range("D47").Value = "0.0005816"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").Value = "0.003513"
$ws.Range("E48").Value = "-2.71%"
$ws.Range("D49").Value = "0.003508"
$ws.Range("E49").Value = "54.23%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.09%"
